# Update the cryptocurrency price/volume table with the latest scraped
# values (GitHub Actions refresh). Price cells (column D) are prefixed
# with a leading apostrophe so Excel stores the numeric-looking text
# (e.g. "1.00", "315.84") as a literal string instead of coercing it to
# a number, matching the sheet's existing text-based formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.737.84"
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("D3").Value = "'2.531.33"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'315.84"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").Value = "'95.69"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").Value = "'0.574"
$ws.Range("E7").Value = "  -0.99%  "
$ws.Range("D9").Value = "'0.532"
$ws.Range("E9").Value = "  -1.21%  "
$ws.Range("D10").Value = "'35.77"
$ws.Range("E10").Value = "  -1.05%  "
$ws.Range("D11").Value = "'0.0806"
$ws.Range("E11").Value = "  -0.23%  "
$ws.Range("D12").Value = "'7.52"
$ws.Range("E12").Value = "  -0.93%  "
$ws.Range("E13").Value = "  -2.25%  "
$ws.Range("D14").Value = "'2.920.44"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "'2.524.15"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "'15.06"
$ws.Range("E16").Value = "  -2.81%  "
$ws.Range("D17").Value = "'0.848"
$ws.Range("E17").Value = "  -1.19%  "
$ws.Range("D18").Value = "'42.805.46"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").Value = "'6.80"
$ws.Range("E19").Value = "  +3.72%  "
$ws.Range("D20").Value = "'12.66"
$ws.Range("E20").Value = "  -2.33%  "
$ws.Range("D21").Value = "'0.0₃0961"
$ws.Range("E21").Value = "  -1.02%  "
$ws.Range("D22").Value = "'69.71"
$ws.Range("E22").Value = "  -2.27%  "
$ws.Range("D23").Value = "'251.30"
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").Value = "'2.95"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").Value = "'2.04"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").Value = "'26.50"
$ws.Range("E26").Value = "  -1.70%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("D29").Value = "'40.40"
$ws.Range("E29").Value = "  +3.91%  "
$ws.Range("D30").Value = "'10.41"
$ws.Range("E30").Value = "  +3.05%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'156.39"
$ws.Range("E32").Value = "  +0.62%  "
$ws.Range("D33").Value = "'2.16"
$ws.Range("E33").Value = "  +3.20%  "
$ws.Range("D34").Value = "'2.72"
$ws.Range("E34").Value = "  +4.10%  "
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "'18.94"
$ws.Range("E36").Value = "  -2.25%  "
$ws.Range("D37").Value = "'0.0779"
$ws.Range("E37").Value = "  -0.97%  "
$ws.Range("D38").Value = "'0.111"
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("E39").Value = "  -1.05%  "
$ws.Range("D40").Value = "'2.29"
$ws.Range("E40").Value = "  +8.06%  "
$ws.Range("D41").Value = "'22.45"
$ws.Range("E41").Value = "  -7.58%  "
$ws.Range("E42").Value = "  -0.85%  "
$ws.Range("B43").Value = "VeChain"
$ws.Range("C43").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D43").Value = "'0.0304"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").Value = "'1.00"
$ws.Range("D45").Value = "'2.026.97"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("E46").Value = "  -3.88%  "
$ws.Range("D47").Value = "'9.08"
$ws.Range("E47").Value = "  +3.03%  "
$ws.Range("D48").Value = "'84.42"
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").Value = "'105.69"
$ws.Range("E49").Value = "  +3.73%  "
$ws.Range("D50").Value = "'74.96"
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("D51").Value = "'2.775.54"
$ws.Range("E51").Value = "  +0.50%  "
